$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 78

$ws.Cells.Item($row, 1).Value = 44306
$ws.Cells.Item($row, 2).Value = 820
$ws.Cells.Item($row, 3).Value = 3406
$ws.Cells.Item($row, 4).Value = 300
$ws.Cells.Item($row, 5).Value = 5259
$ws.Cells.Item($row, 6).Value = 203
$ws.Cells.Item($row, 7).Value = 1816
$ws.Cells.Item($row, 8).Value = 1000
$ws.Cells.Item($row, 9).Value = 28500
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 3517
$ws.Cells.Item($row, 14).Value = 140
$ws.Cells.Item($row, 15).Value = 950
$ws.Cells.Item($row, 16).Value = 42498

# Apply the same date number format used by the rest of column A (e.g. A77)
$ws.Range("A78").NumberFormat = $ws.Range("A77").NumberFormat
